# Applies the "Updated cryptos list" diff: refreshed prices/volumes and a
# one-row insertion (OKB) that shifts rows 9-50 down and drops NEARProtocol.
# All sheet values are plain text (the source data uses t="inlineStr" cells),
# so numeric-looking price strings are entered with a leading apostrophe to
# force Excel to keep them as text instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.818.60'
$ws.Range('E2').Value = '  +3.93%  '
$ws.Range('D3').Value = '1.878.37'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''277.39'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '''0.5288'
$ws.Range('E7').Value = '  +3.88%  '
$ws.Range('D8').Value = '''0.3420'
$ws.Range('E8').Value = '  -3.09%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '''45.13'
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '''0.06947'
$ws.Range('E10').Value = '  +4.19%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').Value = '''20.01'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '''0.8030'
$ws.Range('E12').Value = '  -3.11%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.07738'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.878.64'
$ws.Range('E14').Value = '  +3.46%  '
$ws.Range('D15').Value = '''5.172'
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '''90.13'
$ws.Range('E16').Value = '  +2.96%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '''14.54'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = '''1.000'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '''0.000008040'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.001'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '26.857.75'
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.123.82'
$ws.Range('E22').Value = '  +3.88%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''4.738'
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '''10.02'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '''6.199'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '''2.393'
$ws.Range('E26').Value = '  +8.59%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''146.54'
$ws.Range('E27').Value = '  +3.81%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''1.667'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''17.32'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''113.50'
$ws.Range('E30').Value = '  +3.75%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '''4.337'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''4.298'
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').Value = '''0.08876'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.04889'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''1.174'
$ws.Range('E35').Value = '  +3.09%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''0.7257'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '''2.895'
$ws.Range('E37').Value = '  +1.03%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''3.284'
$ws.Range('E38').Value = '  +4.74%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''2.343'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.01840'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.5123'
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''0.9552'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''116.30'
$ws.Range('E43').Value = '  +4.73%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''6.174'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '''8.077'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = '''1.001'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.4451'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '''0.1339'
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''9.314'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''36.15'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05945'
$ws.Range('E51').Value = '  +1.88%  '
